$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = '37.194.78'
$rng.Style = "Normal"
$ws.Range("E2").Value = '  +0.57%  '
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = '2.059.12'
$rng.Style = "Normal"
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  +0.04%  '
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = '248.66'
$rng.Style = "Normal"
$ws.Range("E5").Value = '  +0.28%  '
$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = '0.666'
$rng.Style = "Normal"
$ws.Range("E6").Value = '  +0.97%  '
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = '58.38'
$rng.Style = "Normal"
$ws.Range("E7").Value = '  +3.94%  '
$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = '1.00'
$rng.Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = '0.385'
$rng.Style = "Normal"
$ws.Range("E9").Value = '  +1.56%  '
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = '0.0789'
$rng.Style = "Normal"
$ws.Range("E10").Value = '  +1.21%  '
$ws.Range("E11").Value = '  +1.78%  '
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = '16.00'
$rng.Style = "Normal"
$ws.Range("E12").Value = '  +1.59%  '
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = '0.917'
$rng.Style = "Normal"
$ws.Range("E13").Value = '  +16.23%  '
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = '2.361.55'
$rng.Style = "Normal"
$ws.Range("E14").Value = '  +1.03%  '
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = '5.82'
$rng.Style = "Normal"
$ws.Range("E15").Value = '  +3.86%  '
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = '2.054.99'
$rng.Style = "Normal"
$ws.Range("E16").Value = '  +0.76%  '
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = '18.94'
$rng.Style = "Normal"
$ws.Range("E17").Value = '  +14.55%  '
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = '37.158.66'
$rng.Style = "Normal"
$ws.Range("E18").Value = '  +0.65%  '
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = '75.48'
$rng.Style = "Normal"
$ws.Range("E19").Value = '  +2.39%  '
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = '0.0₃0909'
$rng.Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = '5.49'
$rng.Style = "Normal"
$ws.Range("E21").Value = '  +3.58%  '
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = '239.65'
$rng.Style = "Normal"
$ws.Range("E22").Value = '  +1.87%  '
$ws.Range("E23").Value = '  -0.07%  '
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = '2.48'
$rng.Style = "Normal"
$ws.Range("E24").Value = '  +4.87%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = '2.21'
$rng.Style = "Normal"
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = '9.60'
$rng.Style = "Normal"
$ws.Range("E26").Value = '  +5.73%  '
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = '171.47'
$rng.Style = "Normal"
$ws.Range("E27").Value = '  +2.35%  '
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = '20.28'
$rng.Style = "Normal"
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = '5.55'
$rng.Style = "Normal"
$ws.Range("E29").Value = '  +18.82%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = '0.125'
$rng.Style = "Normal"
$ws.Range("E30").Value = '  +1.05%  '
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = '1.15'
$rng.Style = "Normal"
$ws.Range("E31").Value = '  +4.14%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = '4.83'
$rng.Style = "Normal"
$ws.Range("E32").Value = '  +9.66%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = '0.0630'
$rng.Style = "Normal"
$ws.Range("E33").Value = '  +3.48%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = '0.0883'
$rng.Style = "Normal"
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = '2.34'
$rng.Style = "Normal"
$ws.Range("E35").Value = '  +6.53%  '
$ws.Range("E36").Value = '  +0.01%  '
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = '1.83'
$rng.Style = "Normal"
$ws.Range("E37").Value = '  +4.75%  '
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = '1.34'
$rng.Style = "Normal"
$ws.Range("E38").Value = '  +0.62%  '
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = '3.13'
$rng.Style = "Normal"
$ws.Range("E39").Value = '  -3.37%  '
$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = '5.18'
$rng.Style = "Normal"
$ws.Range("E40").Value = '  +6.86%  '
$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = '0.102'
$rng.Style = "Normal"
$ws.Range("E41").Value = '  -4.34%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = '0.0227'
$rng.Style = "Normal"
$ws.Range("E42").Value = '  +2.89%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = '101.28'
$rng.Style = "Normal"
$ws.Range("E43").Value = '  +6.40%  '
$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = '1.17'
$rng.Style = "Normal"
$ws.Range("E44").Value = '  +4.87%  '
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = '17.21'
$rng.Style = "Normal"
$ws.Range("E45").Value = '  -0.37%  '
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = '2.43'
$rng.Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = '1.305.75'
$rng.Style = "Normal"
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = '3.82'
$rng.Style = "Normal"
$ws.Range("E48").Value = '  +12.63%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = '2.88'
$rng.Style = "Normal"
$ws.Range("E49").Value = '  +1.49%  '
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = '6.91'
$rng.Style = "Normal"
$ws.Range("E50").Value = '  +3.69%  '
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = '2.249.69'
$rng.Style = "Normal"
$ws.Range("E51").Value = '  +1.07%  '
